# Rebuilds the "Analysis Results" sheet with the updated parameter list:
#  - inserts a new "Date and Time" row at the top
#  - renames/reorders several labels
#  - inserts a new "Cycle Count of battery" row
#  - appends two new "Time spent in 70-80 km/h" / "80-90 km/h" rows
#  - removes the old "Maximum BMS Temperature in C" row
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B1 used to hold the time-duration value/format; it now holds a plain text label,
# so clear its formatting back to General.
$ws.Range("B1").ClearFormats()

$labels = @(
  "Date and Time",
  "Total time taken for the ride",
  "Actual Ampere-hours (Ah)",
  "Actual Watt-hours (Wh)",
  "Starting SoC (Ah)",
  "Ending SoC (Ah)",
  "Starting SoC (%)",
  "Ending SoC (%)",
  "Total distance covered (km)",
  "Total energy consumption(WH/KM)",
  "Total SOC consumed(%)",
  "Mode",
  "Peak Power(kW)",
  "Average Power(kW)",
  "Total Energy Regenerated(kWh)",
  "Regenerative Effectiveness(%)",
  "Highest Cell Voltage(V)",
  "Lowest Cell Voltage(V)",
  "Difference in Cell Voltage(V)",
  "Minimum Temperature(C)",
  "Maximum Temperature(C)",
  "Difference in Temperature(C)",
  "Maximum Fet Temperature-BMS(C)",
  "Maximum Afe Temperature-BMS(C)",
  "Maximum PCB Temperature-BMS(C)",
  "Maximum MCU Temperature(C)",
  "Maximum Motor Temperature(C)",
  "Abnormal Motor Temperature Detected(C)",
  "highest cell temp(C)",
  "lowest cell temp(C)",
  "Difference between Highest and Lowest Cell Temperature at 100% SOC(C)",
  "Battery Voltage(V)",
  "Total energy charged(kWh)",
  "Electricity consumption units(kW)",
  "Cycle Count of battery",
  "Idling time percentage",
  "Time spent in 0-10 km/h",
  "Time spent in 10-20 km/h",
  "Time spent in 20-30 km/h",
  "Time spent in 30-40 km/h",
  "Time spent in 40-50 km/h",
  "Time spent in 50-60 km/h",
  "Time spent in 60-70 km/h",
  "Time spent in 70-80 km/h",
  "Time spent in 80-90 km/h"
)

$bvalues = @(
  "2024-03-12 09:59:28.512000 to 2024-03-12 10:41:56.064000",
  0.02954418981481481,
  32.19162166666667,
  1632.661795212778,
  39.551,
  7.43,
  99,
  18,
  35.8652854330722,
  45.5220633405935,
  81,
  "Custom mode`n76.57%`nEco mode`n19.91%`nSports mode`n0.31%",
  5490.365296,
  -2312.188222960661,
  0.4675920697222222,
  0.02863166099167976,
  3.443,
  3.039,
  0.4039999999999999,
  35,
  47,
  12,
  70,
  68,
  66,
  56,
  97,
  0,
  47,
  35,
  12,
  54,
  1.73834757,
  0.0000001892140771943574049,
  45,
  5.573094799614137,
  3.183372796632465,
  3.288608260983952,
  4.937297202490573,
  24.23923528895905,
  6.585986143997194,
  5.371393492940454,
  18.23204419889503,
  28.43549942997457,
  0
)

for ($i = 0; $i -lt $labels.Length; $i++) {
  $row = $i + 1
  $ws.Cells.Item($row, 1).Value = $labels[$i]
  $ws.Cells.Item($row, 2).Value = $bvalues[$i]
}

# Row 2 ("Total time taken for the ride") keeps the special time-duration format
$ws.Range("B2").NumberFormat = "[hh]:mm:ss"

Write-Host "Done writing" $labels.Length "rows"